$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST")

# Insert a new row at row 11 (shifts existing rows 11..114 down to 12..115)
$ws.Rows.Item(11).Insert()

# Populate the new record inserted at row 11
$ws.Range("A11").Value = "I have severe PTSD from being caught in an avalanche"
$ws.Range("B11").Value = "My partner loves skiing, heli-skiing, and winter mountain climbing."
$ws.Range("C11").Value = "Do you think we should book a backcountry skiing trip in the Alps?"
$ws.Range("D11").Value = "Trauma triggers"

# Match formatting of the surrounding data rows: A:D centered/wrapped, E wrapped only
$dataRange = $ws.Range("A11:D11")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

$ws.Range("E11").WrapText = $true

$ws.Rows.Item(11).RowHeight = 28

# Update the saved selection/active cell to match
$ws.Range("A11:D11").Select()
